$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Step 2 text is shortened, and a new "dashboard" expected result is added
$ws.Range("C3").Value = "Step 2: Login as a service manager"
$ws.Range("D3").Value = "I am redirected to the user's dashboard"

# Row 4: new Step 3 text; expected result keeps same wording (now its own shared string)
$ws.Range("C4").Value = "Step 3:  Go to the create KPI page"
$ws.Range("D4").Value = "I am redirected to the create KPI page"

# Row 5: Step 4 text updated; expected result wording unchanged
$ws.Range("C5").Value = "Step 4: Fill out a KPI for an employee on any of your teams"
$ws.Range("D5").Value = "A new KPI will be added to the database on that employee"

# Row 6: Step 5 text updated; expected result wording updated
$ws.Range("C6").Value = "Step 5: Fill out a KPI for an employee that is not on any of your teams"
$ws.Range("D6").Value = "An error text pops up saying that that access is denied"

# Row 7: brand new Step 6 cells
$ws.Range("C7").Value = "Step 6: Leave Certain fields that are required blank"
$ws.Range("D7").Value = "An error text pops up saying that some of the fields have been left blank"

# Match the updated selection/view shown in the diff
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C7").Select()
